$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.578.39"
$ws.Range("E2").Value = "  +2.10%  "
$ws.Range("D3").Value = "3.391.07"
$ws.Range("E3").Value = "  +2.79%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.74"
$ws.Range("E5").Value = "  +2.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.62"
$ws.Range("E6").Value = "  +7.40%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "3.390.92"
$ws.Range("E8").Value = "  +2.95%  "
$ws.Range("E9").Value = "  +0.82%  "
$ws.Range("E10").Value = "  +1.93%  "
$ws.Range("E11").Value = "  +8.00%  "
$ws.Range("E12").Value = "  +6.36%  "
$ws.Range("D13").Value = "3.969.91"
$ws.Range("E13").Value = "  +3.00%  "
$ws.Range("E15").Value = "  +7.17%  "
$ws.Range("D16").Value = "3.388.35"
$ws.Range("E16").Value = "  +3.09%  "
$ws.Range("E17").Value = "  +4.70%  "
$ws.Range("D18").Value = "61.660.96"
$ws.Range("E18").Value = "  +1.97%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.04"
$ws.Range("E19").Value = "  +5.48%  "
$ws.Range("E20").Value = "  +4.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.38"
$ws.Range("E21").Value = "  +3.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "388.22"
$ws.Range("E22").Value = "  +10.37%  "
$ws.Range("E23").Value = "  +3.43%  "
$ws.Range("D24").Value = "3.527.88"
$ws.Range("E24").Value = "  +3.18%  "
$ws.Range("E25").Value = "  +17.64%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "71.10"
$ws.Range("E27").Value = "  +2.57%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.68"
$ws.Range("E28").Value = "  +6.32%  "
$ws.Range("B29").Value = "Fetch.AI"
$ws.Range("C29").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.59"
$ws.Range("E29").Value = "  +10.93%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.42%  "
$ws.Range("E31").Value = "  +6.32%  "
$ws.Range("E32").Value = "  +5.53%  "
$ws.Range("E33").Value = "  +2.83%  "
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("D35").Value = "3.420.59"
$ws.Range("E35").Value = "  +2.95%  "
$ws.Range("E36").Value = "  +3.50%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.48"
$ws.Range("E37").Value = "  +5.26%  "
$ws.Range("E38").Value = "  +3.18%  "
$ws.Range("E39").Value = "  +4.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "162.32"
$ws.Range("E40").Value = "  +2.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0797"
$ws.Range("E41").Value = "  +6.69%  "
$ws.Range("E42").Value = "  +12.56%  "
$ws.Range("E43").Value = "  +0.19%  "
$ws.Range("E44").Value = "  +7.01%  "
$ws.Range("E45").Value = "  +4.30%  "
$ws.Range("E46").Value = "  +2.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "41.34"
$ws.Range("E47").Value = "  +0.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "24.67"
$ws.Range("E48").Value = "  +8.30%  "
$ws.Range("E49").Value = "  +4.74%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.95"
$ws.Range("E50").Value = "  +7.46%  "
$ws.Range("D51").Value = "2.374.47"
$ws.Range("E51").Value = "  +9.65%  "
